# Update cryptocurrency price (column D) and 1h volume change (column E)
# values to match the latest scrape, row by row, preserving each cell's
# original text (string) representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Prefix with a literal apostrophe so Excel always stores the value
    # as text (never auto-converted to a number), then strip the style
    # Excel applies automatically for "quoted text" back to the sheet
    # default so no unintended formatting is introduced.
    $range = $ws.Range($cell)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextCell "D2" '56.720.82'
$ws.Range("E2").Value = '  +3.29%  '
Set-TextCell "D3" '2.325.14'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  +0.24%  '
Set-TextCell "D5" '520.96'
$ws.Range("E5").Value = '  +3.02%  '
Set-TextCell "D6" '135.21'
$ws.Range("E6").Value = '  +4.44%  '
Set-TextCell "D7" '0.996'
$ws.Range("E7").Value = '  +0.04%  '
Set-TextCell "D8" '0.538'
$ws.Range("E8").Value = '  +1.45%  '
Set-TextCell "D9" '2.349.64'
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("E10").Value = '  +5.77%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  +2.80%  '
$ws.Range("E13").Value = '  +0.68%  '
Set-TextCell "D14" '24.05'
$ws.Range("E14").Value = '  +1.44%  '
Set-TextCell "D15" '2.739.87'
$ws.Range("E15").Value = '  +1.44%  '
Set-TextCell "D16" '56.793.86'
$ws.Range("E16").Value = '  +3.40%  '
$ws.Range("E17").Value = '  +2.04%  '
Set-TextCell "D18" '2.325.46'
$ws.Range("E18").Value = '  +0.43%  '
Set-TextCell "D19" '10.53'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E20").Value = '  +1.48%  '
Set-TextCell "D21" '323.47'
$ws.Range("E21").Value = '  +3.89%  '
Set-TextCell "D22" '6.61'
$ws.Range("E22").Value = '  -0.49%  '
Set-TextCell "D23" '1.00'
$ws.Range("E23").Value = '  +0.24%  '
Set-TextCell "D24" '60.89'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("E25").Value = '  +9.91%  '
Set-TextCell "D26" '0.992'
$ws.Range("E26").Value = '  -0.05%  '
Set-TextCell "D27" '8.00'
$ws.Range("E27").Value = '  +6.57%  '
$ws.Range("E28").Value = '  +13.57%  '
Set-TextCell "D29" '0.0₃0746'
$ws.Range("E29").Value = '  +5.48%  '
$ws.Range("E30").Value = '  +4.92%  '
Set-TextCell "D31" '167.25'
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("E32").Value = '  +1.21%  '
Set-TextCell "D33" '18.39'
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("E34").Value = '  +0.01%  '
Set-TextCell "D35" '0.991'
$ws.Range("E35").Value = '  -0.23%  '
Set-TextCell "D36" '1.26'
$ws.Range("E36").Value = '  +2.40%  '
Set-TextCell "D37" '0.928'
$ws.Range("E37").Value = '  +0.34%  '
Set-TextCell "D38" '4.05'
$ws.Range("E38").Value = '  +4.53%  '
$ws.Range("E39").Value = '  +7.95%  '
Set-TextCell "D40" '37.96'
$ws.Range("E40").Value = '  +3.11%  '
Set-TextCell "D42" '139.55'
$ws.Range("E42").Value = '  +3.84%  '
Set-TextCell "D43" '3.61'
$ws.Range("E43").Value = '  +5.09%  '
Set-TextCell "D44" '5.26'
$ws.Range("E44").Value = '  +7.06%  '
Set-TextCell "D45" '278.96'
$ws.Range("E45").Value = '  +6.83%  '
Set-TextCell "D46" '0.0934'
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("E49").Value = '  +3.25%  '
Set-TextCell "D50" '17.99'
$ws.Range("E50").Value = '  +8.88%  '
$ws.Range("E51").Value = '  +0.88%  '
